$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap contents of columns A and B for rows 1-7 (header + 6 data rows)
for ($r = 1; $r -le 7; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $bVal
    $ws.Cells.Item($r, 2).Value2 = $aVal
}

# Update the selection on the sheet
$ws.Range("B11").Select()
